$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Cells.Item(2, 1).Value = 74
$ws.Cells.Item(2, 2).Value = 133
$ws.Cells.Item(2, 3).Value = 59
$ws.Cells.Item(2, 5).Value = 36
$ws.Cells.Item(2, 6).Value = 38
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(2, 8).Value = 40
$ws.Cells.Item(2, 9).Value = 4

# Update column C values for rows 3-75
$ws.Cells.Item(3, 3).Value = 78
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(5, 3).Value = 32
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(7, 3).Value = 28
$ws.Cells.Item(8, 3).Value = 85
$ws.Cells.Item(9, 3).Value = 105
$ws.Cells.Item(10, 3).Value = 67
$ws.Cells.Item(11, 3).Value = 105
$ws.Cells.Item(12, 3).Value = 45
$ws.Cells.Item(13, 3).Value = 43
$ws.Cells.Item(14, 3).Value = 8
$ws.Cells.Item(15, 3).Value = 71
$ws.Cells.Item(16, 3).Value = 84
$ws.Cells.Item(17, 3).Value = 62
$ws.Cells.Item(18, 3).Value = 121
$ws.Cells.Item(19, 3).Value = 100
$ws.Cells.Item(20, 3).Value = 128
$ws.Cells.Item(21, 3).Value = 21
$ws.Cells.Item(22, 3).Value = 86
$ws.Cells.Item(23, 3).Value = 40
$ws.Cells.Item(24, 3).Value = 34
$ws.Cells.Item(25, 3).Value = 72
$ws.Cells.Item(26, 3).Value = 78
$ws.Cells.Item(27, 3).Value = 75
$ws.Cells.Item(28, 3).Value = 95
$ws.Cells.Item(29, 3).Value = 65
$ws.Cells.Item(30, 3).Value = 79
$ws.Cells.Item(31, 3).Value = 68
$ws.Cells.Item(32, 3).Value = 33
$ws.Cells.Item(33, 3).Value = 129
$ws.Cells.Item(34, 3).Value = 32
$ws.Cells.Item(35, 3).Value = 81
$ws.Cells.Item(36, 3).Value = 17
$ws.Cells.Item(37, 3).Value = 105
$ws.Cells.Item(38, 3).Value = 104
$ws.Cells.Item(39, 3).Value = 48
$ws.Cells.Item(40, 3).Value = 42
$ws.Cells.Item(41, 3).Value = 2
$ws.Cells.Item(42, 3).Value = 3
$ws.Cells.Item(43, 3).Value = 81
$ws.Cells.Item(44, 3).Value = 33
$ws.Cells.Item(45, 3).Value = 109
$ws.Cells.Item(46, 3).Value = 8
$ws.Cells.Item(47, 3).Value = 80
$ws.Cells.Item(48, 3).Value = 35
$ws.Cells.Item(49, 3).Value = 118
$ws.Cells.Item(50, 3).Value = 17
$ws.Cells.Item(51, 3).Value = 24
$ws.Cells.Item(52, 3).Value = 74
$ws.Cells.Item(53, 3).Value = 31
$ws.Cells.Item(54, 3).Value = 104
$ws.Cells.Item(55, 3).Value = 14
$ws.Cells.Item(56, 3).Value = 92
$ws.Cells.Item(57, 3).Value = 19
$ws.Cells.Item(58, 3).Value = 122
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(60, 3).Value = 72
$ws.Cells.Item(61, 3).Value = 122
$ws.Cells.Item(62, 3).Value = 124
$ws.Cells.Item(63, 3).Value = 130
$ws.Cells.Item(64, 3).Value = 61
$ws.Cells.Item(65, 3).Value = 44
$ws.Cells.Item(66, 3).Value = 11
$ws.Cells.Item(67, 3).Value = 117
$ws.Cells.Item(68, 3).Value = 125
$ws.Cells.Item(69, 3).Value = 19
$ws.Cells.Item(70, 3).Value = 108
$ws.Cells.Item(71, 3).Value = 43
$ws.Cells.Item(72, 3).Value = 109
$ws.Cells.Item(73, 3).Value = 56
$ws.Cells.Item(74, 3).Value = 36
$ws.Cells.Item(75, 3).Value = 87
